# "backing up definitions after change"
# Update the definition text for "Grape Price" and move the active
# selection to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grape Price definition (row 15, column B) - reword the definition.
$ws.Range("B15").Value = "Grape price per tonne"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("B16").Select() | Out-Null
